$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    3 = @(0.2917716402565462, 32771568325.09113, 0.7527432677738641, 1133.036916526867, 32771569459.17256)
    4 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    5 = @(0.2917716402565462, 0.04071648406533734, 0.1494219747398047, 0.4942365360607697, 0.9761466351224579)
    6 = @(0.2917716402565462, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.242251378316819)
    7 = @(0.2917716402565462, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 4.630591052383734)
    8 = @(0.1190320826869504, 10.34677158129881, 3.537761648806719, 10.19245300693656, 24.19601831972904)
    9 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 14.05633640148523)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
